$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Sheet1" to "testData"
$ws = $wb.Worksheets.Item(1)
$ws.Name = "testData"

# Set selection to G5 on the active sheet
$ws.Range("G5").Select()

# Set column A width (COM ColumnWidth of 11.1 serialises to OOXML width="12")
$ws.Columns.Item(1).ColumnWidth = 11.1

# Update data in A5 / B5
$ws.Range("A5").Value = "mngr276899"
$ws.Range("B5").Value = "qapydAq"
